$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the "Road" row (Fox / Fox,Wolf) entirely - eliminated terrain ---
$ws.Rows("7:7").Delete() | Out-Null

# --- Rename headers: name -> travelName, previousPath -> previousTerrain ---
$ws.Range("A1").Value = "travelName"
$ws.Range("D1").Value = "previousTerrain"

# --- Rename terrain/travel names in column A (now travel-action phrasing) ---
$ws.Range("A2").Value = "Deeper into the forest"
$ws.Range("A3").Value = "To a meadow nearby"
$ws.Range("A4").Value = "To the mountains"
$ws.Range("A5").Value = "Into a cave in this forest"
$ws.Range("A6").Value = "Follow a stream nearby"
$ws.Range("A7").Value = "To a lake"
$ws.Range("A8").Value = "Into a cave on this mountain"
$ws.Range("A9").Value = "Cross this bridge"
$ws.Range("A10").Value = "In the direction of light forest"

# --- Update previousTerrain (column D) references to match the new names ---
$ws.Range("D2").Value = "In the direction of light forest"
$ws.Range("D5").Value = "Deeper into the forest"
$ws.Range("D8").Value = "To the mountains"
$ws.Range("D9").Value = "Follow a stream nearby"

# --- Column widths (approximate bestFit after the rename) ---
$ws.Columns("A:A").ColumnWidth = 26.6
$ws.Columns("D:D").ColumnWidth = 63.85

# --- Selection as left after the edits ---
$ws.Range("A8").Select() | Out-Null
